$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.299.93"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.575.86"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +2.15%  "
$ws.Range("D5").Value = "'211.42"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").Value = "'46.46"
$ws.Range("E8").Value = "  +5.46%  "
$ws.Range("D9").Value = "'23.81"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'0.0594"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'0.0882"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "1.801.71"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "1.573.63"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "'3.72"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "28.332.29"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'62.63"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'228.11"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "'7.40"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "0.0₃0697"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("D26").Value = "'151.00"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "'15.05"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "'6.50"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").Value = "'3.11"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "1.384.50"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").Value = "'2.57"
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").Value = "'0.800"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").Value = "'1.86"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "'62.59"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "1.712.65"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'85.97"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  +5.11%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0519"
$ws.Range("E51").Value = "  -0.78%  "
